# Automatische test-sync: 2025-08-05 18:21:50
# Adds testmail #8 to the "Logs" sheet, extends the conditional-formatting
# ranges to cover the new row, and refreshes the "Dashboard" pivot-style
# summary counts (Inkoop / Bestellingen now ties Klantenservice / Contact
# at 4, and swaps places in the list).

$wb = $excel.ActiveWorkbook

# ---- Logs sheet: append the new row --------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$newRow = 29
$logs.Range("A" + $newRow).Value = "Kun je nagaan of dit nog leverbaar is?"
$logs.Range("B" + $newRow).Value = "mailmind.test@zohomail.eu"
$logs.Range("C" + $newRow).Value = "Testmail #8: Kun je nagaan of dit nog leverbaar is?"
$logs.Range("D" + $newRow).Value = "Inkoop / Bestellingen"
$logs.Range("E" + $newRow).Value = "Bedankt, we hebben dit doorgestuurd naar inkoop@bedrijf.nl."
$logs.Range("F" + $newRow).Value = "2025-08-05 18:21:41"
$logs.Range("G" + $newRow).Value = "Ja"
$logs.Range("H" + $newRow).Value = "Ja"
$logs.Range("I" + $newRow).Value = "Nee"
$logs.Range("J" + $newRow).Value = "Nee"

# ---- Extend conditional formatting ranges to include the new row ---------
foreach ($col in @("D", "G", "H", "I", "J")) {
    $fcs = $logs.Range($col + "2").FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fc = $fcs.Item($i)
        $fc.ModifyAppliesToRange($logs.Range($col + "2:" + $col + $newRow))
    }
}

# ---- Dashboard sheet: update summary counts -------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A3").Value = "Inkoop / Bestellingen"
$dash.Range("B3").Value = 4
$dash.Range("A4").Value = "Klantenservice / Contact"
$dash.Range("B4").Value = 4
